# Fixed issue where SummarySection.summary_variable included variables for
# all sections: reorder the three report blocks on sheet "Page one" so
# that "Last name begins with C" now comes first, then "Birth to last
# vote for David", then "Birth to last vote" (previously "Birth to last
# vote" was first, "Last name begins with C" second, "Birth to last vote
# for David" last).

$xlPasteFormats = -4122
$xlPasteValues = -4163

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Page one")

function Copy-Block($srcRange, $dstRange) {
    $ws.Range($srcRange).Copy()
    $ws.Range($dstRange).PasteSpecial($xlPasteFormats)
    $ws.Range($srcRange).Copy()
    $ws.Range($dstRange).PasteSpecial($xlPasteValues)
    $excel.CutCopyMode = $false
}

# Places a merged section title: copies the format from a template
# merged title cell (staged in a scratch row, since the original title
# rows get wiped by Clear()), (re)creates the merge without leaving the
# "split border" look that Range.Merge() introduces, and sets the title
# text without disturbing the merge.
function Set-Title($titleCell, $mergedRange, $text) {
    if (-not $ws.Range($titleCell).MergeCells) {
        $ws.Range($mergedRange).Merge()
    }
    $ws.Range("B100:D100").Copy()
    $ws.Range($mergedRange).PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = $false
    $ws.Range($titleCell).Value = $text
}

# --- Stage the three existing blocks (values + formats) far away from
# --- the working area so copying them back in their new order never
# --- overwrites source data that hasn't been copied yet. Title/
# --- column-header rows only use columns B:D (column A must stay
# --- untouched/empty there); data rows use the full A:D range. Title
# --- text is staged separately (plain text, no merge needed there).

$titleBirthToLastVote = $ws.Range("B1").Value2
$titleLastNameC = $ws.Range("B8").Value2
$titleDavid = $ws.Range("B13").Value2

# Template of the title-row format (bold/fill/border), staged so it
# survives the Range("A1:D15").Clear() call further below.
Copy-Block "B1:D1" "B100:D100"

# Block A: "Birth to last vote"            rows 1-5   (5 rows)
Copy-Block "B2:D2"   "B102:D102"
Copy-Block "A3:D5"   "A103:D105"
# Block B: "Last name begins with C"       rows 8-10  (3 rows)
Copy-Block "B9:D9"   "B112:D112"
Copy-Block "A10:D10" "A113:D113"
# Block C: "Birth to last vote for David"  rows 13-15 (3 rows)
Copy-Block "B14:D14" "B122:D122"
Copy-Block "A15:D15" "A123:D123"

# --- Unmerge the old title cells and clear the original block area
# --- (row 1's merge is reused as-is for the new row-1 title, so leave
# --- it merged).
$ws.Range("B8:D8").UnMerge()
$ws.Range("B13:D13").UnMerge()
$ws.Range("A1:D15").Clear()

# --- Copy the staged blocks back into their new locations / order.

# New rows 1-3: "Last name begins with C" (staged block B)
Set-Title "B1" "B1:D1" $titleLastNameC
Copy-Block "B112:D112" "B2:D2"
Copy-Block "A113:D113" "A3:D3"

# New rows 6-8: "Birth to last vote for David" (staged block C)
Set-Title "B6" "B6:D6" $titleDavid
Copy-Block "B122:D122" "B7:D7"
Copy-Block "A123:D123" "A8:D8"

# New rows 11-15: "Birth to last vote" (staged block A)
Set-Title "B11" "B11:D11" $titleBirthToLastVote
Copy-Block "B102:D102" "B12:D12"
Copy-Block "A103:D105" "A13:D15"

# --- Clean up the staging area.
$ws.Range("A100:D123").Clear()

$ws.Range("A1").Select()
